$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:F1) ---
$ws.Range("A1").Value = "fullname"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "sid"
$ws.Range("D1").Value = "classs"
$ws.Range("E1").Value = "roll"
$ws.Range("F1").Value = "barcode"

# "roll" header gets wrap text formatting
$ws.Range("E1").WrapText = $true

# --- Row 2 : AHETESHAM URRAB (existing person, columns re-mapped) ---
$ws.Range("A2").Value = "AHETESHAM URRAB"
$ws.Range("B2").Value = "ahete@gmail.com"
$ws.Range("C2").Value = 3803556
$ws.Range("D2").Value = "TYBBACA"
$ws.Range("E2").Value = 6623

# --- Row 3 : FIRDOUS (new person) ---
$ws.Range("A3").Value = "FIRDOUS"
$ws.Range("B3").Value = "firdous@gmail.com"
$ws.Range("C3").Value = 4087370
$ws.Range("D3").Value = "TYBBACA"
$ws.Range("E3").Value = 6656

# --- Row 4 : ALID SHAUKH (text overwritten in-place; hyperlink rel left as-is) ---
$ws.Range("A4").Value = "ALID SHAUKH"
$ws.Range("B4").Value = "alid@gmail.com"
$ws.Range("C4").Value = 3800563
$ws.Range("D4").Value = "TYBBACA"
$ws.Range("E4").Value = 6674

# --- Row 5 : MOHAMMED MAAZ SHAIKH (new person) ---
$ws.Range("A5").Value = "MOHAMMED MAAZ SHAIKH"
$ws.Range("B5").Value = "maaz@gmail.com"
$ws.Range("C5").Value = 3801330
$ws.Range("D5").Value = "TYBBACA"
$ws.Range("E5").Value = 6675

# --- Row 6 : SALAUDDIN KHAN (new person) ---
$ws.Range("A6").Value = "SALAUDDIN KHAN"
$ws.Range("B6").Value = "salaudin@gmail.com"
$ws.Range("C6").Value = 3801322
$ws.Range("D6").Value = "TYBBACA"
$ws.Range("E6").Value = 6617

# --- Row 7 : ABUSHAMA (new person) ---
$ws.Range("A7").Value = "ABUSHAMA"
$ws.Range("B7").Value = "abu@gmail.com"
$ws.Range("C7").Value = 3801332
$ws.Range("D7").Value = "TYBBACA"
$ws.Range("E7").Value = 6640

# --- Row 8 : UZMA SAYYED (new person) ---
$ws.Range("A8").Value = "UZMA SAYYED"
$ws.Range("B8").Value = "uzma@gmail.com"
$ws.Range("C8").Value = 3803389
$ws.Range("D8").Value = "TYBBACA"
$ws.Range("E8").Value = 6626

# --- Hyperlinks for the new rows' email column (B2:B4 already had links; leave untouched) ---
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:maaz@gmail.com")
$ws.Range("B5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:salaudin@gmail.com")
$ws.Range("B6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:abu@gmail.com")
$ws.Range("B7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:uzma@gmail.com")
$ws.Range("B8").Style = "Hyperlink"

# --- column widths ---
$ws.Columns("A").ColumnWidth = 23.7109375
$ws.Columns("F").ColumnWidth = 12.140625

# --- final selection ---
$ws.Range("E8").Select()
